$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(2,4,8,9,12,15,17,18,21,24,30,34,35,36,37,39,43,46,47,49,50,53,56,58,65,68,72,74,77,80,84,85,86,93,97,99,100,107,110,113,114,116,120,123,129,135,141,145,146,152,156,162,163,165,167,168,169,172,173,176,180,182,188,193)

foreach ($r in $rows) {
    $ws.Range("G$r`:N$r").Value2 = "N/A"
}
